$d = $word.ActiveDocument

# --- 1) Bump the "Heading 2" style font size by one point (11pt -> 12pt, i.e.
#        w:sz 22 -> 24 half-points) so level-2 headings render slightly bigger.
$heading2 = $d.Styles("Heading 2")
$heading2.Font.Size = 12

# --- 2) Relocate the hidden "_GoBack" bookmark: it currently sits mid-sentence
#        in the "When you click the Knit..." paragraph (an artifact of the last
#        edit position); move it to the very start of the document (before the
#        title text), which is where Word drops it after a fresh edit there.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Collapsed Range(0,0) passed straight to Bookmarks.Add snaps to cover the
# first word instead of staying zero-length, so nudge into place indirectly:
# insert a throwaway character, bookmark right after it, then remove the
# character again - the bookmark collapses back to a true zero-length mark
# at position 0.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$afterTemp = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterTemp)
$d.Range(0, 1).Delete()

# --- 3) The paragraph that used to hold the "_GoBack" bookmark had been split
#        into two runs around it ("... will be" / " generated ..."); now that
#        the bookmark is gone, fold the text back into a single run.
$d.Content.Find.Execute(
    " button a document will be generated that includes both content as well as the output of any embedded R code chunks within the document. You can embed an R code chunk like this:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " button a document will be generated that includes both content as well as the output of any embedded R code chunks within the document. You can embed an R code chunk like this:",
    2) | Out-Null
